$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row for "Update Voice Bot Agent Status ( PUT )" API, appended after
# the existing last row (45) at row 48 (rows 46/47 are intentionally left
# blank, matching the spacing pattern already used by rows 36/39/42/45).

$aText = @'
Update Voice Bot  Agent Status ( PUT )

'@

$bText = 'http://1msg.1point1.in:3001/api/auth/j-v1/update/agent/status/'

$cText = @'
{
    "user_id" : 10,
    "agent_id" : 78,
    "agent_status" : true
}
'@

$dText = @'
curl --location --request PUT 'http://1msg.1point1.in:3001/api/auth/j-v1/update/agent/status/' \
--header 'Content-Type: application/json' \
--data '{
    "user_id" : 10,
    "agent_id" : 78,
    "agent_status" : true
}'

'@

# Write cell values in column order (A, B, C, D) so the new shared-string
# entries are appended in the same order as the target workbook.
$ws.Range("A48").Value = $aText
$ws.Range("B48").Value = $bText
$ws.Range("C48").Value = $cText
$ws.Range("D48").Value = $dText

# A48 / C48 / D48 use the wrap-text style already used elsewhere in the
# sheet (e.g. C45/D45); copy that formatting across without disturbing the
# row height. B48 keeps the default (no special) style, same as this row
# in the target file.
$ws.Range("C45").Copy() | Out-Null
$ws.Range("A48").PasteSpecial(-4122) | Out-Null
$ws.Range("C48").PasteSpecial(-4122) | Out-Null
$ws.Range("D48").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

# Match the row's rendered height from the source workbook.
$ws.Rows.Item(48).RowHeight = 151.8

# Update the active selection to D48 (matches the authored workbook).
$ws.Range("D48").Select() | Out-Null
